$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.889.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.047.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "525.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.67%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +4.98%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.60"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.15%  "
$ws.Range("E10").Value = "  +8.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.370"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.11%  "
$ws.Range("E12").Value = "  +2.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.568.50"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000172"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +17.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "57.777.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.21"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.041.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "341.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.500"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.174"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0₃0973"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.59%  "
$ws.Range("E27").Value = "  +0.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.23%  "
$ws.Range("E30").Value = "  +7.44%  "
$ws.Range("E31").Value = "  +4.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "156.75"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.15%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.35"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.72%  "
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.93"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +11.90%  "
$ws.Range("E38").Value = "  +2.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.076.30"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.74"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.88"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.85%  "
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.55%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.662"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.49%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.322.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.42%  "
$ws.Range("E46").Value = "  +2.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0248"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0895"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.59%  "
